$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 72; this shifts the existing rows 72-130 down to 73-131
# (matching the row-level diff, which is a pure insert-and-shift of a new
# weekly record ahead of the former row 72).
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new record's data.
$ws.Cells.Item(72, 1).Value = 10
$ws.Cells.Item(72, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(72, 3).Value = 'La Araucanía'
$ws.Cells.Item(72, 4).Value = 44827
$ws.Cells.Item(72, 5).Value = 9
$ws.Cells.Item(72, 6).Value = 100114002
$ws.Cells.Item(72, 7).Value = 'Camote'
$ws.Cells.Item(72, 8).Value = 'Sin especificar'
$ws.Cells.Item(72, 9).Value = 'Primera'
$ws.Cells.Item(72, 10).Value = 40
$ws.Cells.Item(72, 11).Value = 20000
$ws.Cells.Item(72, 12).Value = 20000
$ws.Cells.Item(72, 13).Value = 20000
$ws.Cells.Item(72, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(72, 15).Value = 'Perú'
$ws.Cells.Item(72, 16).Value = 1000
$ws.Cells.Item(72, 17).Value = 20
$ws.Cells.Item(72, 18).Value = 'Hortaliza'
